$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("current punchlist")

# --- Row 6: "Color map points..." task is now DONE, released in 1.0.0 on 2022-05-16 ---
$ws.Range("B6").Value = "DONE"
$ws.Range("C6").Value = 44696
$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("E6").Value = "Color map points by alert status - combo of level and trend."
$ws.Range("D6").Value = "1.0.0"

# --- Row 9: mark as POSTPONED ---
$ws.Range("B9").Value = "POSTPONED"

# --- Row 12: mark as WORKING ---
$ws.Range("B12").Value = "WORKING"

# --- Row 14: mark as WORKING, add a note, and grow the row height for the extra note ---
$ws.Range("B14").Value = "WORKING"
$ws.Range("F14").Value = "Tried numerous times but doesn't seem to work as advertised."
$ws.Rows("14").RowHeight = 44

# --- Row 16: mark as WORKING ---
$ws.Range("B16").Value = "WORKING"

# --- Row 17: mark as PENDING ---
$ws.Range("B17").Value = "PENDING"

# --- Row 18: "Deploy from github directly." task is now DONE, released in 0.9.1 on 5/1/2022 ---
$ws.Range("B18").Value = "DONE"
$ws.Range("C18").Value = 44681
$ws.Range("C11").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D18").Value = "0.9.1"

# --- Column B is wider now to fit status values ---
$ws.Columns("B").ColumnWidth = 22

# --- Update selection to reflect where editing left off ---
$ws.Range("D19").Select()
$excel.CutCopyMode = $false
